$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "qntVinculos" column (old column C) - shifts tipoTabela..taxa
# one column to the left (D->C, E->D, F->E, G->F, H->G).
$ws.Columns.Item(3).Delete()

# ----- Header row -----
$ws.Range("A1").Value = "nomeVendedor"
$ws.Range("B1").Value = "cnpj"
$ws.Range("C1").Value = "tipoTabela"
$ws.Range("D1").Value = "Origem"
$ws.Range("E1").Value = "seTaxa"
$ws.Range("F1").Value = "tipoTaxa"
$ws.Range("G1").Value = "taxa"

# ----- Row 2 -----
$ws.Range("A2").Value = "Erik"

$ws.Range("B2").NumberFormat = "general"
$ws.Range("B2").Value = 20400219000195
$ws.Range("B2").NumberFormat = "@"

$ws.Range("C2").Value = "Package 4"
$ws.Range("D2").Value = "SSA"
$ws.Range("E2").Value = "N"

$ws.Range("F2").NumberFormat = "general"
$ws.Range("F2").Value = 0

$ws.Range("G2").NumberFormat = "general"
$ws.Range("G2").Value = 0

# ----- Row 3 -----
$ws.Range("A3").Value = "Erik"

$ws.Range("B3").NumberFormat = "general"
$ws.Range("B3").Value = 20400219000195
$ws.Range("B3").NumberFormat = "@"

$ws.Range("C3").Value = ".COM 2"
$ws.Range("D3").Value = "SSA"
$ws.Range("E3").Value = "N"

$ws.Range("F3").NumberFormat = "general"
$ws.Range("F3").Value = 0

$ws.Range("G3").NumberFormat = "general"
$ws.Range("G3").Value = 0

# ----- Row 4 - cleared out, keep its (text) style -----
$ws.Range("A4:G4").ClearContents()

# ----- Column widths -----
$ws.Columns.Item(1).ColumnWidth = 15.140625
$ws.Columns.Item(3).ColumnWidth = 10.42578125
$ws.Columns.Item(5).ColumnWidth = 10.5703125

# ----- Page setup -----
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ----- Selection -----
$ws.Range("I3").Select()
